# Applies the Golem_Profits market-data refresh captured in the commit diff.
# Each entry updates columns H-N (price/profit columns) for one leve row on one sheet.
$wb = $excel.ActiveWorkbook

$rowUpdates = @(
    @{ Sheet = "ALC"; Row = 40; Sets = @{ H=1599; I=1599; J=0; K=1599; L=0; M=-1424 }; Clears = @("N") }
    @{ Sheet = "ALC"; Row = 45; Sets = @{ H=16000; I=0; J=16000; K=0; L=48000; N=-48384 }; Clears = @("M") }
    @{ Sheet = "ALC"; Row = 53; Sets = @{ H=168.83333; I=174.6; J=140; K=174.6; L=140; M=462.4; N=-1414 }; Clears = @() }
    @{ Sheet = "ALC"; Row = 99; Sets = @{ H=1002.8182; I=600; J=1338.5; K=1800; L=4015.5; M=-302; N=-7011.5 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 12; Sets = @{ H=0; J=0; L=0 }; Clears = @("N") }
    @{ Sheet = "ARM"; Row = 38; Sets = @{ H=4966666.5; I=0; J=4966666.5; K=0; L=4966666.5; N=-4967600.5 }; Clears = @("M") }
    @{ Sheet = "ARM"; Row = 45; Sets = @{ H=3544.25; I=1479.1428; K=1479.1428; M=-1102.1428 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 61; Sets = @{ H=2542; I=2038.8; J=3800; K=2038.8; L=3800; M=-1826.8; N=-4224 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 74; Sets = @{ H=336039.16; I=402487.1; K=402487.1; M=-401613.1 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 77; Sets = @{ H=336039.16; I=402487.1; K=2012435.5; M=-2008067.5 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 97; Sets = @{ H=4966.5; I=4983; J=4950; K=4983; L=4950; M=-4487; N=-5942 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 101; Sets = @{ H=25333.334; J=25333.334; L=25333.334; N=-31823.334 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 109; Sets = @{ H=70000; J=70000; L=70000; N=-72774 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 122; Sets = @{ H=8937.538; I=8811.75; K=26435.25; M=-23985.25 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 124; Sets = @{ H=40499; J=71000; L=71000; N=-80820 }; Clears = @() }
    @{ Sheet = "ARM"; Row = 136; Sets = @{ H=2542; I=2038.8; J=3800; K=6116.4; L=11400; M=-3566.4; N=-16500 }; Clears = @() }
    @{ Sheet = "BSM"; Row = 33; Sets = @{ H=80000; I=0; K=0 }; Clears = @("M") }
    @{ Sheet = "BSM"; Row = 94; Sets = @{ H=2288.4614; I=1916.6666; J=2607.1428; K=1916.6666; L=2607.1428; M=-1465.6666; N=-3509.1428 }; Clears = @() }
    @{ Sheet = "BSM"; Row = 105; Sets = @{ H=748.3333; I=748.3333; K=748.3333; M=998.6667 }; Clears = @() }
    @{ Sheet = "BSM"; Row = 134; Sets = @{ H=3083.125; I=3083.125; K=9249.375; M=-6714.375 }; Clears = @() }
    @{ Sheet = "CRP"; Row = 86; Sets = @{ H=2813.8; J=2250; L=2250; N=-4496 }; Clears = @() }
    @{ Sheet = "CRP"; Row = 89; Sets = @{ H=2813.8; J=2250; L=11250; N=-22482 }; Clears = @() }
    @{ Sheet = "CRP"; Row = 132; Sets = @{ H=0; J=0; L=0 }; Clears = @("N") }
    @{ Sheet = "CRP"; Row = 134; Sets = @{ H=921.75; I=921.75; K=2765.25; M=-230.25 }; Clears = @() }
    @{ Sheet = "CUL"; Row = 11; Sets = @{ H=112.77778; I=136; J=31.5; K=408; L=94.5; M=-268; N=-374.5 }; Clears = @() }
    @{ Sheet = "GSM"; Row = 74; Sets = @{ H=52999; J=52999; L=52999; N=-54871 }; Clears = @() }
    @{ Sheet = "GSM"; Row = 77; Sets = @{ H=52999; J=52999; L=158997; N=-168357 }; Clears = @() }
    @{ Sheet = "GSM"; Row = 97; Sets = @{ H=1001999.7; I=2999.5; K=2999.5; M=-2503.5 }; Clears = @() }
    @{ Sheet = "GSM"; Row = 120; Sets = @{ H=0; J=0; L=0 }; Clears = @("N") }
    @{ Sheet = "LTW"; Row = 22; Sets = @{ H=2119.8; I=1766.3334; J=2650; K=1766.3334; L=2650; M=-1471.3334; N=-3240 }; Clears = @() }
    @{ Sheet = "LTW"; Row = 27; Sets = @{ H=2119.8; I=1766.3334; J=2650; K=1766.3334; L=2650; M=-1659.3334; N=-2864 }; Clears = @() }
    @{ Sheet = "LTW"; Row = 46; Sets = @{ H=0; J=0; L=0 }; Clears = @("N") }
    @{ Sheet = "LTW"; Row = 122; Sets = @{ H=3362.9092; I=3099; K=9297; M=-6847 }; Clears = @() }
    @{ Sheet = "LTW"; Row = 124; Sets = @{ H=0; I=0; K=0 }; Clears = @("M") }
    @{ Sheet = "LTW"; Row = 136; Sets = @{ H=724143.5600000001; I=844334.2; J=3000; K=2533002.6; L=9000; M=-2530452.6; N=-14100 }; Clears = @() }
    @{ Sheet = "WVR"; Row = 54; Sets = @{ H=44498.5; I=50000; J=38997; K=50000; L=38997; M=-49480; N=-40037 }; Clears = @() }
    @{ Sheet = "WVR"; Row = 62; Sets = @{ H=2000; J=2000; L=2000; N=-3248 }; Clears = @() }
    @{ Sheet = "WVR"; Row = 65; Sets = @{ H=2000; J=2000; L=10000; N=-16240 }; Clears = @() }
    @{ Sheet = "WVR"; Row = 122; Sets = @{ H=810.25; I=873.5; K=2620.5; M=-170.5 }; Clears = @() }
    @{ Sheet = "WVR"; Row = 124; Sets = @{ H=0; J=0; L=0 }; Clears = @("N") }
    @{ Sheet = "WVR"; Row = 136; Sets = @{ H=2749.5; I=2749.5; K=8248.5; M=-5698.5 }; Clears = @() }
)

foreach ($u in $rowUpdates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    foreach ($col in $u.Sets.Keys) {
        $ws.Range("$col$($u.Row)").Value = $u.Sets[$col]
    }
    foreach ($col in $u.Clears) {
        $ws.Range("$col$($u.Row)").Value = $null
    }
}

Write-Host "Updated $($rowUpdates.Count) leve rows with refreshed market data."